{"js": "// Correction de faute : suppression du mot \"etc\" en trop dans la phrase\n// \"Partir des bases (rappels, paradigmes, introduction aux technologies etc \u2026)\"\n// qui devient\n// \"Partir des bases (rappels, paradigmes, introduction aux technologies\u2026)\"\nconst body = context.document.body;\n\nconst results = body.search(\" etc \", { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Correction de faute : suppression du mot \"etc\" en trop dans la phrase\n# \"Partir des bases (rappels, paradigmes, introduction aux technologies etc \u2026)\"\n# qui devient\n# \"Partir des bases (rappels, paradigmes, introduction aux technologies\u2026)\"\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \" etc \"\n$find.Replacement.Text = \"\"\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $wdReplaceAll) | Out-Null\n"}
